$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the run_id column) so we append correctly
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row

# New log rows to append
$newRows = @(
    @{ run_id = 79; rss_url_id = 1; date = "2024-06-16 14:10:57"; response = 200; item_count = 9 },
    @{ run_id = 80; rss_url_id = 2; date = "2024-06-16 14:10:58"; response = 200; item_count = 1 }
)

foreach ($entry in $newRows) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $entry.run_id
    $ws.Cells.Item($lastRow, 2).Value = $entry.rss_url_id
    $ws.Cells.Item($lastRow, 3).Value = $entry.date
    $ws.Cells.Item($lastRow, 4).Value = $entry.response
    $ws.Cells.Item($lastRow, 5).Value = $entry.item_count
}
